# Reorder the workbook's sheet tabs: "总计" (summary) should come before
# "2020-Q4" so the tab order becomes [总计, 2020-Q4] instead of
# [2020-Q4, 总计]. No cell data is changed - the two worksheets simply
# swap tab positions (matches commit: "update data with resort sheetname").

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")

# Move "总计" so it sits immediately before the current first worksheet,
# i.e. it becomes the first tab; "2020-Q4" then naturally follows as
# the second tab.
$totalSheet.Move($wb.Worksheets.Item(1))
